$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.109.01"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "'1.601.75"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'301.93"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("D7").Value = "'0.3783"
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("D8").Value = "'0.3653"
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("D9").Value = "'49.83"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "'1.269"
$ws.Range("E10").Value = "  -6.06%  "
$ws.Range("D11").Value = "'0.08160"
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "'6.596"
$ws.Range("E14").Value = "  -6.59%  "
$ws.Range("D15").Value = "'0.00001261"
$ws.Range("E15").Value = "  -4.20%  "
$ws.Range("D16").Value = "'7.417"
$ws.Range("E16").Value = "  -8.47%  "
$ws.Range("D17").Value = "'1.599.02"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").Value = "'91.79"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("D19").Value = "'0.06843"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "'18.51"
$ws.Range("E20").Value = "  -5.84%  "
$ws.Range("D21").Value = "'6.591"
$ws.Range("E21").Value = "  -5.50%  "
$ws.Range("D22").Value = "'0.5561"
$ws.Range("E22").Value = "  -6.89%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'13.04"
$ws.Range("E24").Value = "  -5.42%  "
$ws.Range("D25").Value = "'23.117.26"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").Value = "'2.729"
$ws.Range("E27").Value = "  -8.17%  "
$ws.Range("D28").Value = "'21.18"
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("D29").Value = "'150.35"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "'5.263"
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("D31").Value = "'132.80"
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("D32").Value = "'2.390"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").Value = "'6.876"
$ws.Range("E33").Value = "  -12.55%  "
$ws.Range("D34").Value = "'1.774.59"
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("D35").Value = "'0.9623"
$ws.Range("E35").Value = "  -5.88%  "
$ws.Range("D36").Value = "'0.07717"
$ws.Range("E36").Value = "  -5.97%  "
$ws.Range("D37").Value = "'6.296"
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("D38").Value = "'0.02733"
$ws.Range("E38").Value = "  -6.19%  "
$ws.Range("D39").Value = "'0.2560"
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("D40").Value = "'0.08897"
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").Value = "'10.09"
$ws.Range("E41").Value = "  -7.13%  "
$ws.Range("D42").Value = "'1.371"
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").Value = "'0.7091"
$ws.Range("E43").Value = "  -6.42%  "
$ws.Range("D44").Value = "'12.65"
$ws.Range("E44").Value = "  -7.56%  "
$ws.Range("D45").Value = "'15.40"
$ws.Range("E45").Value = "  -6.84%  "
$ws.Range("D46").Value = "'0.6622"
$ws.Range("E46").Value = "  -4.84%  "
$ws.Range("E47").Value = "  -5.86%  "
$ws.Range("D48").Value = "'0.9995"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'3.996"
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("D50").Value = "'132.70"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'0.07940"
$ws.Range("E51").Value = "  -4.23%  "
